# "Write excel datas add" - populate the "Puan" (Points) column (E) on the
# active sheet: a header label plus the point totals for the rows that have
# scores recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Puan"
$ws.Range("E2").Value = 85
$ws.Range("E5").Value = 78
$ws.Range("E10").Value = 72
